# Update cryptocurrency price/volume data per the diff.
# D-column values that look numeric are prefixed with a leading apostrophe
# so Excel stores them as literal text (matching the source inlineStr cells)
# instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $ws.Range("D2").Value = "'56.663.36"
    $ws.Range("E2").Value = "  +0.21%  "
    $ws.Range("D3").Value = "'2.386.65"
    $ws.Range("E3").Value = "  +0.63%  "
    $ws.Range("E4").Value = "  -0.10%  "
    $ws.Range("D5").Value = "'504.87"
    $ws.Range("E5").Value = "  +0.82%  "
    $ws.Range("D6").Value = "'132.54"
    $ws.Range("E6").Value = "  +2.98%  "
    $ws.Range("E7").Value = "  -0.03%  "
    $ws.Range("D9").Value = "'2.390.80"
    $ws.Range("E9").Value = "  -0.08%  "
    $ws.Range("D10").Value = "'0.0975"
    $ws.Range("E10").Value = "  +2.24%  "
    $ws.Range("E11").Value = "  +0.88%  "
    $ws.Range("E12").Value = "  +2.88%  "
    $ws.Range("D13").Value = "'4.67"
    $ws.Range("E13").Value = "  +1.89%  "
    $ws.Range("D14").Value = "'2.813.15"
    $ws.Range("E14").Value = "  -0.06%  "
    $ws.Range("D15").Value = "'56.595.26"
    $ws.Range("E15").Value = "  -0.79%  "
    $ws.Range("D16").Value = "'21.68"
    $ws.Range("E16").Value = "  +0.90%  "
    $ws.Range("D18").Value = "'2.376.94"
    $ws.Range("E18").Value = "  -4.26%  "
    $ws.Range("D19").Value = "'10.19"
    $ws.Range("E19").Value = "  +0.88%  "
    $ws.Range("E20").Value = "  +0.96%  "
    $ws.Range("D21").Value = "'309.80"
    $ws.Range("E21").Value = "  -0.14%  "
    $ws.Range("D22").Value = "'6.26"
    $ws.Range("E22").Value = "  +1.19%  "
    $ws.Range("E23").Value = "  -0.08%  "
    $ws.Range("D24").Value = "'5.58"
    $ws.Range("E24").Value = "  -4.61%  "
    $ws.Range("D25").Value = "'66.26"
    $ws.Range("E25").Value = "  +1.38%  "
    $ws.Range("E26").Value = "  -0.86%  "
    $ws.Range("E27").Value = "  +0.24%  "
    $ws.Range("E28").Value = "  -0.06%  "
    $ws.Range("D29").Value = "'7.35"
    $ws.Range("E29").Value = "  +2.36%  "
    $ws.Range("D30").Value = "'175.67"
    $ws.Range("E30").Value = "  +1.11%  "
    $ws.Range("D31").Value = "'0.0₃0726"
    $ws.Range("E31").Value = "  +2.65%  "
    $ws.Range("E32").Value = "  -0.09%  "
    $ws.Range("E33").Value = "  +2.86%  "
    $ws.Range("E34").Value = "  -3.65%  "
    $ws.Range("E35").Value = "  +0.04%  "
    $ws.Range("D36").Value = "'0.997"
    $ws.Range("E36").Value = "  +0.11%  "
    $ws.Range("D37").Value = "'17.81"
    $ws.Range("E37").Value = "  +0.57%  "
    $ws.Range("E38").Value = "  -0.74%  "
    $ws.Range("E39").Value = "  +1.82%  "
    $ws.Range("D40").Value = "'36.80"
    $ws.Range("E40").Value = "  +2.78%  "
    $ws.Range("D41").Value = "'0.819"
    $ws.Range("E41").Value = "  +6.59%  "
    $ws.Range("E42").Value = "  +1.37%  "
    $ws.Range("D43").Value = "'132.21"
    $ws.Range("E43").Value = "  +2.71%  "
    $ws.Range("E44").Value = "  +1.25%  "
    $ws.Range("E45").Value = "  +1.03%  "
    $ws.Range("E46").Value = "  -0.44%  "
    $ws.Range("E47").Value = "  +1.14%  "
    $ws.Range("D48").Value = "'246.54"
    $ws.Range("E48").Value = "  -2.58%  "
    $ws.Range("D49").Value = "'0.0484"
    $ws.Range("E49").Value = "  +0.46%  "
    $ws.Range("E50").Value = "  +1.87%  "
    $ws.Range("D51").Value = "'17.17"
    $ws.Range("E51").Value = "  +7.64%  "
